# Weekly data refresh: prepend the newest "Zapallo italiano" price record
# (Macroferia Regional de Talca) as the new first row of the data block,
# pushing the existing historical rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a fresh blank row just above the current first data row (row 521),
# shifting rows 521:596 down to 522:597.
$ws.Rows("521:521").Insert()

# Populate the newly inserted row 521 with the new week's record.
$ws.Range("A521").Value = 5
$ws.Range("B521").Value = "Macroferia Regional de Talca"
$ws.Range("C521").Value = "Maule"
$ws.Range("D521").Value = 45127
$ws.Range("E521").Value = 7
$ws.Range("F521").Value = 100112032
$ws.Range("G521").Value = "Zapallo italiano"
$ws.Range("H521").Value = "Sin especificar"
$ws.Range("I521").Value = "Primera"
$ws.Range("J521").Value = 300
$ws.Range("K521").Value = 14000
$ws.Range("L521").Value = 14000
$ws.Range("M521").Value = 14000
$ws.Range("N521").Value = "`$/caja 50 unidades"
$ws.Range("O521").Value = "Región del Maule"
$ws.Range("P521").Value = 280
$ws.Range("Q521").Value = 50
$ws.Range("R521").Value = "Hortaliza"
